# 10Th - MB for single stock and added new group
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column B. This shifts the existing
# B,C,D,E columns (Jun_17,Jun_15,Jun_13,Jun_10 + analyst data) over to E,F,G,H
$ws.Range("B:D").Insert()

# New header cells for the 3 inserted date columns
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Give the new C and D columns (and reaffirm F,G,H) the same custom
# width as the rest of the data columns in the sheet
$ws.Columns("C:H").ColumnWidth = 7.14

# Fill the newly inserted columns (B,C,D) with "UN" for every existing
# analyst data row (rows 2-27), matching the existing B/C/D "UN" pattern.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Add two new rows for a new analyst/group section at the bottom
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
